$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue $ws "D2" "64.449.53"
Set-TextValue $ws "E2" "  +2.17%  "
Set-TextValue $ws "D3" "2.676.00"
Set-TextValue $ws "E3" "  +2.95%  "
Set-TextValue $ws "E4" "  +0.06%  "
Set-TextValue $ws "D5" "596.43"
Set-TextValue $ws "E5" "  +2.22%  "
Set-TextValue $ws "D6" "148.15"
Set-TextValue $ws "E6" "  +0.18%  "
Set-TextValue $ws "E7" "  +0.05%  "
Set-TextValue $ws "D8" "0.593"
Set-TextValue $ws "E8" "  -0.96%  "
Set-TextValue $ws "E9" "  +0.43%  "
Set-TextValue $ws "D10" "5.67"
Set-TextValue $ws "E10" "  -0.16%  "
Set-TextValue $ws "E11" "  +0.02%  "
Set-TextValue $ws "D12" "0.359"
Set-TextValue $ws "E12" "  +1.23%  "
Set-TextValue $ws "D13" "28.00"
Set-TextValue $ws "E13" "  +2.72%  "
Set-TextValue $ws "D14" "3.159.26"
Set-TextValue $ws "E14" "  +3.03%  "
Set-TextValue $ws "D15" "64.371.99"
Set-TextValue $ws "E15" "  +2.21%  "
Set-TextValue $ws "E16" "  +0.49%  "
Set-TextValue $ws "D17" "2.683.02"
Set-TextValue $ws "E17" "  +3.24%  "
Set-TextValue $ws "D18" "11.44"
Set-TextValue $ws "E18" "  +0.72%  "
Set-TextValue $ws "D19" "346.78"
Set-TextValue $ws "E19" "  +1.02%  "
Set-TextValue $ws "D20" "4.42"
Set-TextValue $ws "E20" "  +0.21%  "
Set-TextValue $ws "D21" "6.90"
Set-TextValue $ws "E21" "  +1.78%  "
Set-TextValue $ws "E22" "  +0.20%  "
Set-TextValue $ws "D23" "68.90"
Set-TextValue $ws "E23" "  +2.60%  "
Set-TextValue $ws "B24" "SuiNetwork"
Set-TextValue $ws "C24" "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue $ws "D24" "1.62"
Set-TextValue $ws "E24" "  +11.09%  "
Set-TextValue $ws "B25" "Fetch.AI"
Set-TextValue $ws "C25" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws "D25" "1.67"
Set-TextValue $ws "E25" "  +4.51%  "
Set-TextValue $ws "E26" "  -1.05%  "
Set-TextValue $ws "E27" "  +1.73%  "
Set-TextValue $ws "D28" "8.02"
Set-TextValue $ws "E28" "  +1.49%  "
Set-TextValue $ws "D29" "0.999"
Set-TextValue $ws "E29" "  -0.10%  "
Set-TextValue $ws "D30" "531.41"
Set-TextValue $ws "E30" "  +14.07%  "
Set-TextValue $ws "E31" "  +3.59%  "
Set-TextValue $ws "E32" "  +11.43%  "
Set-TextValue $ws "D33" "0.0₃0828"
Set-TextValue $ws "E33" "  +0.69%  "
Set-TextValue $ws "D34" "175.50"
Set-TextValue $ws "E35" "  +0.07%  "
Set-TextValue $ws "D36" "0.404"
Set-TextValue $ws "E36" "  +0.43%  "
Set-TextValue $ws "D37" "19.34"
Set-TextValue $ws "E37" "  +0.61%  "
Set-TextValue $ws "D38" "4.68"
Set-TextValue $ws "E38" "  +1.90%  "
Set-TextValue $ws "E39" "  +4.07%  "
Set-TextValue $ws "D40" "173.21"
Set-TextValue $ws "E40" "  +8.77%  "
Set-TextValue $ws "D41" "0.999"
Set-TextValue $ws "E41" "  +0.04%  "
Set-TextValue $ws "D42" "40.72"
Set-TextValue $ws "E42" "  +3.03%  "
Set-TextValue $ws "D43" "3.79"
Set-TextValue $ws "E43" "  -0.16%  "
Set-TextValue $ws "D44" "21.92"
Set-TextValue $ws "E44" "  +4.12%  "
Set-TextValue $ws "E45" "  -0.42%  "
Set-TextValue $ws "D46" "0.0552"
Set-TextValue $ws "E46" "  +0.89%  "
Set-TextValue $ws "E47" "  +1.70%  "
Set-TextValue $ws "D48" "0.0965"
Set-TextValue $ws "E48" "  -0.96%  "
Set-TextValue $ws "D49" "18.86"
Set-TextValue $ws "E49" "  +1.49%  "
Set-TextValue $ws "D50" "1.78"
Set-TextValue $ws "E50" "  +2.91%  "
Set-TextValue $ws "E51" "  -0.46%  "
